# Relabel several header cells in row 1 of the Monthly Report template.
#   B1: "Site Type"              -> "Entry Point"
#   G1: "Tester Name"            -> "Lab Manager Name"
#   J1: "Tester Contact Number"  -> "Lab Manager Contact Number"
# and leave the sheet scrolled/selected on the (now relabeled) J1 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Entry Point"
$ws.Range("G1").Value = "Lab Manager Name"
$ws.Range("J1").Value = "Lab Manager Contact Number"

# Scroll the view back to the top-left (A1) and leave the active
# selection on J1, matching the saved sheet view in the workbook.
$ws.Range("J1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
